$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HomePage")

$ws.Activate()

$ws.Range("A6").Value = "lbl_search_result"
$ws.Range("B6").Value = "XPath"
$ws.Range("C6").Value = "//ul[@class='product_list grid row']//a[contains(text(),'searchText')]"

$ws.Range("A6").Select()
